# Coronavirus.xlsx — "Add files via upload" edit
#
# Summary of the change (per the OOXML diff):
#   * Nurses sheet   : the footnoted "Brazil*" / "Chile*" country labels are
#                       replaced by the plain "Brazil" / "Chile" labels.
#   * Medicos sheet  : the three extra breakdown columns (Generalist /
#                       Specialist / Medical doctors not further defined)
#                       are removed, shrinking the table from A1:G6 to
#                       A1:D6. The surviving columns/values are untouched.
#   * Workbook       : the "data__1" defined name (scoped to Medicos) is
#                       re-pointed at the now-smaller A1:D6 range.
#   * Selection       : the Medicos sheet's remembered cell selection moves
#                       to A11.
#
# Shared strings that become unused ("Brazil*", "Chile*", the three
# breakdown-column headers) are dropped automatically by the engine when it
# re-serialises xl/sharedStrings.xml, so nothing special needs to be done
# for that here.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Nurses sheet: "Brazil*"/"Chile*" -> "Brazil"/"Chile"
# ---------------------------------------------------------------------
$nurses = $wb.Worksheets.Item("Nurses")
$nurses.Range("A3").Value = "Brazil"
$nurses.Range("A4").Value = "Chile"

# ---------------------------------------------------------------------
# 2. Medicos sheet: drop the Generalist/Specialist/Not-further-defined
#    columns (E:G), leaving Country/Year/Medical doctors per 10k/Medical
#    doctors number (A:D).
# ---------------------------------------------------------------------
$medicos = $wb.Worksheets.Item("Medicos")
$medicos.Range("E1:G6").Delete() | Out-Null

# ---------------------------------------------------------------------
# 3. Workbook-level defined name now covers the smaller range.
# ---------------------------------------------------------------------
$wb.Names.Item("data__1").RefersTo = "=Medicos!`$A`$1:`$D`$6"

# ---------------------------------------------------------------------
# 4. Remembered selection: Medicos -> A11. Re-select the Nurses sheet
#    afterwards so it stays the active tab (as in the source file).
# ---------------------------------------------------------------------
$medicos.Range("A11").Select() | Out-Null
$nurses.Range("A10").Select() | Out-Null
